# Add the example kid row to the (previously empty) data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "sex"
$ws.Range("C1").Value = "age (year.month)"
$ws.Range("D1").Value = "two_consistent_pairing"
$ws.Range("E1").Value = "fam_order"
$ws.Range("F1").Value = "fam_grouping"
$ws.Range("G1").Value = "shape_order"
$ws.Range("H1").Value = "test_condition"
$ws.Range("I1").Value = "test_trial_order"
$ws.Range("J1").Value = "pretest_choice_most_like"
$ws.Range("K1").Value = "pretest_choice_least_like"
$ws.Range("L1").Value = "test_choice_most_like"
$ws.Range("M1").Value = "test_choice_least_like"

# ---- Example kid data (row 2) ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "M"
$ws.Range("C2").Value = 3.8
$ws.Range("D2").Value = "H"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "A"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "H"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 2

# ---- Cosmetics: widen the data columns to fit their header text, like Excel's
# "best fit" would after typing the headers in, then leave the selection where
# the author left off (age column, data row). ----
$ws.Range("C1:M2").EntireColumn.AutoFit() | Out-Null

$ws.Range("C2").Select() | Out-Null
